$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1146.3889
$ws.Range("I19").Value = 1144.6
$ws.Range("J19").Value = 1148.625
$ws.Range("K19").Value = 1144.6
$ws.Range("L19").Value = 1148.625
$ws.Range("M19").Value = -969.5999999999999
$ws.Range("N19").Value = -1498.625

$ws.Range("H62").Value = 7642.2173
$ws.Range("I62").Value = 6936.85
$ws.Range("K62").Value = 6936.85
$ws.Range("M62").Value = -6312.85

$ws.Range("H63").Value = 109995
$ws.Range("J63").Value = 109995
$ws.Range("L63").Value = 109995
$ws.Range("N63").Value = -111243

$ws.Range("H65").Value = 7642.2173
$ws.Range("I65").Value = 6936.85
$ws.Range("K65").Value = 34684.25
$ws.Range("M65").Value = -31564.25

$ws.Range("H66").Value = 109995
$ws.Range("J66").Value = 109995
$ws.Range("L66").Value = 329985
$ws.Range("N66").Value = -336225

$ws.Range("H95").Value = 97659.664
$ws.Range("J95").Value = 97659.664
$ws.Range("L95").Value = 97659.664
$ws.Range("N95").Value = -103151.664

$ws.Range("H105").Value = 49330
$ws.Range("J105").Value = 49330
$ws.Range("L105").Value = 49330
$ws.Range("N105").Value = -56318

$ws.Range("H109").Value = 99995
$ws.Range("J109").Value = 99995
$ws.Range("L109").Value = 99995
$ws.Range("N109").Value = -102769

$ws.Range("H114").Value = 99994.5
$ws.Range("J114").Value = 99994.5
$ws.Range("L114").Value = 99994.5
$ws.Range("N114").Value = -108672.5

$ws.Range("H120").Value = 116999.5
$ws.Range("J120").Value = 116999.5
$ws.Range("L120").Value = 116999.5
$ws.Range("N120").Value = -126675.5

$ws.Range("H124").Value = 54329.668
$ws.Range("J124").Value = 54329.668
$ws.Range("L124").Value = 54329.668
$ws.Range("N124").Value = -64149.668

$ws.Range("H126").Value = 139986
$ws.Range("J126").Value = 139986
$ws.Range("L126").Value = 139986
$ws.Range("N126").Value = -149866

$ws.Range("H130").Value = 116998.5
$ws.Range("J130").Value = 116998.5
$ws.Range("L130").Value = 116998.5
$ws.Range("N130").Value = -127038.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17938230
$ws.Range("I32").Value = 18264034
$ws.Range("K32").Value = 18264034
$ws.Range("M32").Value = -18263747

$ws.Range("H63").Value = 4430
$ws.Range("J63").Value = 4912.5
$ws.Range("L63").Value = 4912.5
$ws.Range("N63").Value = -6284.5

$ws.Range("H66").Value = 4430
$ws.Range("J66").Value = 4912.5
$ws.Range("L66").Value = 24562.5
$ws.Range("N66").Value = -31426.5

$ws.Range("H74").Value = 3080.4614
$ws.Range("I74").Value = 3129.3333
$ws.Range("K74").Value = 3129.3333
$ws.Range("M74").Value = -2255.3333

$ws.Range("H77").Value = 3080.4614
$ws.Range("I77").Value = 3129.3333
$ws.Range("K77").Value = 15646.6665
$ws.Range("M77").Value = -11278.6665

$ws.Range("H80").Value = 99997.836
$ws.Range("J80").Value = 99997.836
$ws.Range("L80").Value = 99997.836
$ws.Range("N80").Value = -101993.836

$ws.Range("H83").Value = 99997.836
$ws.Range("J83").Value = 99997.836
$ws.Range("L83").Value = 299993.508
$ws.Range("N83").Value = -309977.508

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 65651
$ws.Range("I20").Value = 85751.336
$ws.Range("J20").Value = 5350
$ws.Range("K20").Value = 85751.336
$ws.Range("L20").Value = 5350
$ws.Range("M20").Value = -85504.336
$ws.Range("N20").Value = -5844

$ws.Range("H107").Value = 2738.8096
$ws.Range("I107").Value = 2971.8235
$ws.Range("J107").Value = 1748.5
$ws.Range("K107").Value = 2971.8235
$ws.Range("L107").Value = 1748.5
$ws.Range("M107").Value = -1051.8235
$ws.Range("N107").Value = -5588.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7637
$ws.Range("I31").Value = 4455.5
$ws.Range("J31").Value = 14000
$ws.Range("K31").Value = 4455.5
$ws.Range("L31").Value = 14000
$ws.Range("M31").Value = -4160.5
$ws.Range("N31").Value = -14590

$ws.Range("H34").Value = 7637
$ws.Range("I34").Value = 4455.5
$ws.Range("J34").Value = 14000
$ws.Range("K34").Value = 4455.5
$ws.Range("L34").Value = 14000
$ws.Range("M34").Value = -4253.5
$ws.Range("N34").Value = -14404

$ws.Range("H68").Value = 49979.7
$ws.Range("J68").Value = 49979.7
$ws.Range("L68").Value = 49979.7
$ws.Range("N68").Value = -51477.7

$ws.Range("H71").Value = 49979.7
$ws.Range("J71").Value = 49979.7
$ws.Range("L71").Value = 149939.1
$ws.Range("N71").Value = -157427.1

$ws.Range("H99").Value = 3103.111
$ws.Range("I99").Value = 2847.5715
$ws.Range("K99").Value = 2847.5715
$ws.Range("M99").Value = -1349.5715

$ws.Range("H126").Value = 3103.111
$ws.Range("I126").Value = 2847.5715
$ws.Range("K126").Value = 8542.7145
$ws.Range("M126").Value = -6072.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7950.2856
$ws.Range("I3").Value = 7950.2856
$ws.Range("K3").Value = 23850.8568
$ws.Range("M3").Value = -23738.8568

$ws.Range("H18").Value = 2842.2856
$ws.Range("I18").Value = 832.3333
$ws.Range("K18").Value = 2496.9999
$ws.Range("M18").Value = -2327.9999

$ws.Range("H38").Value = 41.15
$ws.Range("J38").Value = 66.44444
$ws.Range("L38").Value = 199.33332
$ws.Range("N38").Value = -893.33332

$ws.Range("H122").Value = 436.25
$ws.Range("I122").Value = 436.25
$ws.Range("K122").Value = 3926.25
$ws.Range("M122").Value = -1476.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 509.85715
$ws.Range("I2").Value = 497
$ws.Range("K2").Value = 497
$ws.Range("M2").Value = -384

$ws.Range("H102").Value = 1619.7368
$ws.Range("I102").Value = 1698.2858
$ws.Range("K102").Value = 1698.2858
$ws.Range("M102").Value = -76.28580000000011

$ws.Range("H107").Value = 1009.8182
$ws.Range("I107").Value = 1085.8
$ws.Range("K107").Value = 1085.8
$ws.Range("M107").Value = 834.2

$ws.Range("H126").Value = 2694.7778
$ws.Range("I126").Value = 2549.25
$ws.Range("J126").Value = 2811.2
$ws.Range("K126").Value = 7647.75
$ws.Range("L126").Value = 8433.599999999999
$ws.Range("M126").Value = -5177.75
$ws.Range("N126").Value = -13373.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()

$ws.Range("H40").Value = 3907.3684
$ws.Range("J40").Value = 3917.0715
$ws.Range("L40").Value = 3917.0715
$ws.Range("N40").Value = -4189.0715

$ws.Range("H61").Value = 3248.6155
$ws.Range("I61").Value = 1368.3077
$ws.Range("J61").Value = 5128.923
$ws.Range("K61").Value = 1368.3077
$ws.Range("L61").Value = 5128.923
$ws.Range("M61").Value = -1166.3077
$ws.Range("N61").Value = -5532.923

$ws.Range("H113").Value = 3248.6155
$ws.Range("I113").Value = 1368.3077
$ws.Range("J113").Value = 5128.923
$ws.Range("K113").Value = 1368.3077
$ws.Range("L113").Value = 5128.923
$ws.Range("M113").Value = 801.6922999999999
$ws.Range("N113").Value = -9468.922999999999

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H133").Value = 28180
$ws.Range("J133").Value = 28180
$ws.Range("L133").Value = 28180
$ws.Range("N133").Value = -33240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -41108

$ws.Range("H62").Value = 5435.2856
$ws.Range("I62").Value = 3002
$ws.Range("J62").Value = 5840.8335
$ws.Range("K62").Value = 3002
$ws.Range("L62").Value = 5840.8335
$ws.Range("M62").Value = -2378
$ws.Range("N62").Value = -7088.8335

$ws.Range("H65").Value = 5435.2856
$ws.Range("I65").Value = 3002
$ws.Range("J65").Value = 5840.8335
$ws.Range("K65").Value = 15010
$ws.Range("L65").Value = 29204.1675
$ws.Range("M65").Value = -11890
$ws.Range("N65").Value = -35444.1675

$ws.Range("H100").Value = 672.5
$ws.Range("I100").Value = 240.33333
$ws.Range("K100").Value = 480.66666
$ws.Range("M100").Value = 60.33334000000002

$ws.Range("H108").Value = 115000
$ws.Range("J108").Value = 115000
$ws.Range("L108").Value = 115000
$ws.Range("N108").Value = -122680

$ws.Range("H126").Value = 1523.0714
$ws.Range("I126").Value = 1580.8182
$ws.Range("K126").Value = 4742.4546
$ws.Range("M126").Value = -2272.4546

$ws.Range("H132").Value = 4073.5833
$ws.Range("J132").Value = 3985.5
$ws.Range("L132").Value = 11956.5
$ws.Range("N132").Value = -17016.5
